# Apply updated crypto price/volume figures per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.084.53'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '3.519.05'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('D5').Value = "'572.73"
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = "'183.66"
$ws.Range('E6').Value = '  -2.60%  '
$ws.Range('E7').Value = '  -2.35%  '
$ws.Range('D8').Value = '3.512.19'
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = "'0.187"
$ws.Range('E10').Value = '  +4.89%  '
$ws.Range('D11').Value = "'0.640"
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('D12').Value = "'53.95"
$ws.Range('E12').Value = '  -3.80%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '4.087.60'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('D17').Value = '3.517.57'
$ws.Range('E17').Value = '  -1.81%  '
$ws.Range('D18').Value = '69.076.35'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = "'12.53"
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').Value = "'539.62"
$ws.Range('E21').Value = '  +13.87%  '
$ws.Range('E22').Value = '  -0.82%  '
$ws.Range('D23').Value = "'20.71"
$ws.Range('E23').Value = '  +9.44%  '
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('D25').Value = "'4.44"
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').Value = "'94.76"
$ws.Range('E26').Value = '  +6.74%  '
$ws.Range('D27').Value = "'10.94"
$ws.Range('E27').Value = '  -1.04%  '
$ws.Range('D28').Value = "'2.92"
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').Value = "'9.17"
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').Value = "'31.55"
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('D31').Value = "'7.24"
$ws.Range('E31').Value = '  -5.99%  '
$ws.Range('D32').Value = "'12.70"
$ws.Range('E32').Value = '  +4.85%  '
$ws.Range('D33').Value = "'64.25"
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('E34').Value = '  -4.59%  '
$ws.Range('D35').Value = "'574.65"
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').Value = "'38.00"
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = "'3.05"
$ws.Range('E39').Value = '  +5.80%  '
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('D42').Value = "'3.10"
$ws.Range('E42').Value = '  -3.70%  '
$ws.Range('E43').Value = '  -4.38%  '
$ws.Range('D44').Value = "'3.56"
$ws.Range('E44').Value = '  +7.03%  '
$ws.Range('E45').Value = '  -4.65%  '
$ws.Range('D46').Value = '3.211.32'
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range('D47').Value = "'0.0441"
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = "'9.13"
$ws.Range('E48').Value = '  -4.74%  '
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').Value = "'135.89"
$ws.Range('E51').Value = '  -0.96%  '
